# Implement first version of lot sizing rules:
#  - Generic sheet: NrBuckets (B4) goes from 4 to 5
#  - ForecastedAverageDemand and ForcastedStandardDeviation sheets: add a new
#    bucket row (row 6) mirroring the pattern of the existing rows.

$wb = $excel.ActiveWorkbook

# --- Generic sheet: bump NrBuckets from 4 to 5 ---
$genericSheet = $wb.Worksheets.Item("Generic")
$genericSheet.Range("B4").Value = 5

# Row 6 values shared by both forecast sheets (A..N)
$newRowValues = @(4, 0, 0, 0, 0, 0, 1, 1, 1, 1, 0, 0, 0, 0)

$sheetNames = @("ForecastedAverageDemand", "ForcastedStandardDeviation")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    for ($col = 1; $col -le $newRowValues.Length; $col++) {
        $ws.Cells.Item(6, $col).Value = $newRowValues[$col - 1]
    }

    # Match the styling of column A used in the preceding bucket rows (A2:A5)
    $ws.Range("A5").Copy()
    $ws.Range("A6").PasteSpecial(-4122)
}
